$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Año")

# New header for column B: "Fecha"
$ws.Range("B1").Value = "Fecha"

# First date (1/1/2017), written as its serial number so no stray time-of-day
# fraction is introduced, then formatted with the built-in short-date format
# (numFmtId 14).
$ws.Range("B2").Value = 42736
$ws.Range("B2").NumberFormat = "mm-dd-yy"

# Copy the formatting (not the value) from B2 down to B3:B5 so that all the
# date cells share the very same style/number-format entry.
$ws.Range("B2").Copy()
$ws.Range("B3:B5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Now fill in the remaining date serials (1/1/2018, 1/1/2019, 1/1/2020).
$ws.Range("B3").Value = 43101
$ws.Range("B4").Value = 43466
$ws.Range("B5").Value = 43831

# Restore the selection state recorded in the saved workbook.
$ws.Range("C5:C6").Select()
